$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''70.733.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.78%  '

$ws.Range("D3").Value = '''3.571.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.22%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''582.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.07%  '

$ws.Range("D6").Value = '''186.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.85%  '

$ws.Range("E7").Value = '  +2.38%  '

$ws.Range("D8").Value = '''3.560.36'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.16%  '

$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("E10").Value = '  +21.37%  '

$ws.Range("D11").Value = '''0.651'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.22%  '

$ws.Range("D12").Value = '''54.52'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.12%  '

$ws.Range("D13").Value = '''0.0000319'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.22%  '

$ws.Range("E14").Value = '  +0.62%  '

$ws.Range("D15").Value = '''4.137.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.06%  '

$ws.Range("D16").Value = '''70.774.73'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.91%  '

$ws.Range("D17").Value = '''19.25'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.12%  '

$ws.Range("D18").Value = '''3.570.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.25%  '

$ws.Range("D19").Value = '''12.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.30%  '

$ws.Range("D20").Value = '''573.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.17%  '

$ws.Range("E21").Value = '  +0.53%  '

$ws.Range("E22").Value = '  -0.85%  '

$ws.Range("E23").Value = '  -8.66%  '

$ws.Range("D24").Value = '''4.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.88%  '

$ws.Range("D25").Value = '''4.90'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.07%  '

$ws.Range("D26").Value = '''94.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.15%  '

$ws.Range("D27").Value = '''11.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.37%  '

$ws.Range("E28").Value = '  +2.02%  '

$ws.Range("D29").Value = '''9.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.53%  '

$ws.Range("D30").Value = '''32.41'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.95%  '

$ws.Range("D31").Value = '''7.20'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.48%  '

$ws.Range("E32").Value = '  -1.73%  '

$ws.Range("D33").Value = '''0.116'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.63%  '

$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = '''3.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +15.27%  '

$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").Value = '''63.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.32%  '

$ws.Range("D36").Value = '''548.68'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.33%  '

$ws.Range("E37").Value = '  +13.68%  '

$ws.Range("E38").Value = '  +4.86%  '

$ws.Range("B39").Value = 'InjectiveProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D39").Value = '''38.22'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.30%  '

$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '''0.0₃0805'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.59%  '

$ws.Range("D41").Value = '''0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.07%  '

$ws.Range("D42").Value = '''3.581.95'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.13%  '

$ws.Range("E43").Value = '  +4.69%  '

$ws.Range("E44").Value = '  +3.61%  '

$ws.Range("D45").Value = '''0.0463'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.15%  '

$ws.Range("D46").Value = '''3.48'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.01%  '

$ws.Range("E47").Value = '  -1.27%  '

$ws.Range("D48").Value = '''9.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.16%  '

$ws.Range("E49").Value = '  +2.46%  '

$ws.Range("E50").Value = '  +14.73%  '

$ws.Range("D51").Value = '''0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.12%  '
